# Add 6 new LeetCode tracking rows (240-245) to the "刷题记录表" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=240; A=3432; B="Count Partitions with Even Sum Difference"; C="#array #math "; D="easy";   E=0.5; F=0.5; G=5;  H="2025-12-05"; I="2025-12-05" },
    @{ Row=241; A=198;  B="House Robber"; C="#array #dynamic-programming  "; D="medium"; E=0.5; F=0.5; G=8;  H="2025-12-05"; I="2025-12-05" },
    @{ Row=242; A=152;  B="Maximum Product Subarray"; C="#array #dynamic-programming "; D="medium"; E=0; F=1; G=50; H="2025-12-05"; I="2025-12-05" },
    @{ Row=243; A=3578; B="Count Partitions With Max-Min Difference at Most K"; C="#sliding-window #dynamic-programming "; D="medium"; E=0; F=0; G=40; H="2025-12-06"; I="2025-12-06" },
    @{ Row=244; A=1523; B="Count Odd Numbers in an Interval Range"; C="#math"; D="easy"; E=1; F=0; G=3; H="2025-12-07"; I="2025-12-07" },
    @{ Row=245; A=1925; B="Count Square Sum Triples"; C="#math"; D="easy"; E=1; F=0; G=3; H="2025-12-08"; I="2025-12-08" }
)

$heights = @{ 240=34; 241=34; 242=34; 243=51; 244=34; 245=34 }

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Rows.Item($r.Row).RowHeight = $heights[$r.Row]
}

$excel.ActiveWindow.TopLeftCell = "A241"
$ws.Range("A248").Select()
